$d = $word.ActiveDocument

function Replace-Text($oldText, $newText) {
    $rng = $d.Content
    $result = $rng.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
    return $result
}

# --- Edit 1: Remove the leading two sentences from the 2nd Introduction paragraph ---
# (they get relocated to the end of the 1st Introduction paragraph below)
$old1 = 'The increased variability of weather caused by climate change has a negative impact on food supply, access to food, and food quality. The dry season has been longer than it has ever been in recent years. It starts earlier and ends later in the year. Drought becomes more likely as a result of this. This risk has a direct impact on the amount of water available for agriculture and home consumption, as scarcity will drive up prices. Farmers with already modest incomes will find it increasingly difficult to secure water for their crops as costs rise. '
$found1 = Replace-Text $old1 ""
if (-not $found1) { Write-Output "EDIT1 (trim para4 lead) NOT FOUND" }

# --- Edit 2: Replace the 1st Introduction paragraph with the new opening + relocated sentences ---
$old2 = 'Climate change refers to changes in the atmosphere that are driven by both natural and man-made variables such as the earth''s orbit, volcanic activity, and crustal motions. Artificial causes such as the growth in greenhouse gas and aerosol concentrations climate change caused by global warming, which refers to the average increase in global temperature, has emerged as a megatrend that will result in huge future world changes. Climate change''s effects can also be seen in other ways around the world, such as rising sea levels, melting glaciers, northward relocation of plant habitats, changes in animal habitats, rising ocean temperatures, shorter winters, and early arrival of spring. '
$new2 = 'The Philippines'' vulnerability to the effects of global climate change has long been known. The country came in seventh place out of more than 180 countries in the Global Climate Risk Index. Over the last 20 years, severe weather-related disasters have had the greatest impact (Kreft & Eckstein, 2014). According to the United Nations University Institute for Environment and Human Security''s World Risk Index Report, it ranked third in terms of climate change vulnerability, particularly exposure to natural climate disasters (2012). Extreme weather events are predicted to become more common as a result of climate change. The increased variability of weather caused by climate change has a negative impact on food supply, access to food, and food quality. The dry season has been longer than it has ever been in recent years. It starts earlier and ends later in the year. Drought becomes more likely as a result of this. This risk has a direct impact on the amount of water available for agriculture and home consumption, as scarcity will drive up prices. Farmers with already modest incomes will find it increasingly difficult to secure water for their crops as costs rise. '
$found2 = Replace-Text $old2 $new2
if (-not $found2) { Write-Output "EDIT2 (intro para1 rewrite) NOT FOUND" }

# --- Edit 3: Rewrite the Significance-of-the-project paragraph ---
$old3 = 'The goal of this analysis is to determine how climate change affects agriculture in the Philippines by applying exploratory data analysis to temperature and weather changes, as well as crop production value. The researchers will be able to uncover strategies to avoid or mitigate climate change in the Philippines by constructing exploratory data analysis for this study. Farmers will profit as well, since they will be able to adjust to the unexpected change in weather. And to the economy, it will result higher income for the people, proper price of the products and the reducing poverty.    '
$new3 = 'The goal of this analysis is to examine how climate change affects farmers in the Philippines by applying exploratory data analysis to climate change indicators, as well as Crops Statistics of the Philippines 2016-2020, and Rice and Corn Stocks Inventory, January 2022. The researchers will be able to uncover strategies to avoid or mitigate climate change in the Philippines by constructing exploratory data analysis for this study. Farmers will profit as well, since they will be able to adjust to the unexpected change in weather. '
$found3 = Replace-Text $old3 $new3
if (-not $found3) { Write-Output "EDIT3 (significance rewrite) NOT FOUND" }

Write-Output "DONE"
